$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Valor Mora" total
$ws.Range("E11").Value = 77818

# Update worker / period counts
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1

# Update salary basico for the first worker row (VIVIANA)
$ws.Range("G16").Value = 1423500

# Remove the four intermediate rows (DEIBIS, and the 3 extra JULIA periods).
# This shifts the BEATRIZ row (old row 21) up to become row 17, bringing its
# own data and formatting with it, and shifts the signature block up too.
$ws.Rows("17:20").Delete()
